$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 211.4614666666667
$ws.Range("H2").Value = 634.3844
$ws.Range("I2").Value = 0.2421062275331183
$ws.Range("J2").Value = 0.2421062275331183
$ws.Range("M2").Value = 12.19250666666667
$ws.Range("N2").Value = 36.57752
$ws.Range("O2").Value = 0.2503750524168468
$ws.Range("P2").Value = 0.2503750524168468
$ws.Range("Q2").Value = 2578.245342076445
$ws.Range("R2").Value = 23204.208078688
$ws.Range("S2").Value = 0.06061735940904953
$ws.Range("T2").Value = 0.06061735940904953
$ws.Range("G3").Value = 211.4614666666667
$ws.Range("H3").Value = 634.3844
$ws.Range("I3").Value = 0.2421062275331183
$ws.Range("J3").Value = 0.2421062275331183
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.03913333333333333
$ws.Range("N3").Value = 0.1174
$ws.Range("O3").Value = 0.0008036091882046079
$ws.Range("P3").Value = 0.000803609188204608
$ws.Range("Q3").Value = 8.275192062222223
$ws.Range("R3").Value = 74.47672856000001
$ws.Range("S3").Value = 0.0001945587889671693
$ws.Range("T3").Value = 0.0001945587889671693
$ws.Range("G4").Value = 211.4614666666667
$ws.Range("H4").Value = 634.3844
$ws.Range("I4").Value = 0.2421062275331183
$ws.Range("J4").Value = 0.2421062275331183
$ws.Range("M4").Value = 23.13337333333333
$ws.Range("N4").Value = 69.40011999999999
$ws.Range("O4").Value = 0.4750474795102416
$ws.Range("P4").Value = 0.4750474795102416
$ws.Range("Q4").Value = 4891.817054014221
$ws.Range("R4").Value = 44026.35348612799
$ws.Range("S4").Value = 0.1150119531633409
$ws.Range("T4").Value = 0.1150119531633409
$ws.Range("G5").Value = 211.4614666666667
$ws.Range("H5").Value = 634.3844
$ws.Range("I5").Value = 0.2421062275331183
$ws.Range("J5").Value = 0.2421062275331183
$ws.Range("M5").Value = 13.33195766666667
$ws.Range("N5").Value = 39.995873
$ws.Range("O5").Value = 0.273773858884707
$ws.Range("P5").Value = 0.273773858884707
$ws.Range("Q5").Value = 2819.195321731245
$ws.Range("R5").Value = 25372.7578955812
$ws.Range("S5").Value = 0.0662823561717607
$ws.Range("T5").Value = 0.06628235617176069
$ws.Range("I6").Value = 0.08842543241393927
$ws.Range("J6").Value = 0.08842543241393927
$ws.Range("M6").Value = 12.19250666666667
$ws.Range("N6").Value = 36.57752
$ws.Range("O6").Value = 0.2503750524168468
$ws.Range("P6").Value = 0.2503750524168468
$ws.Range("Q6").Value = 941.6629285636533
$ws.Range("R6").Value = 8474.966357072881
$ws.Range("S6").Value = 0.02213952227562238
$ws.Range("T6").Value = 0.02213952227562238
$ws.Range("I7").Value = 0.08842543241393927
$ws.Range("J7").Value = 0.08842543241393927
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.03913333333333333
$ws.Range("N7").Value = 0.1174
$ws.Range("O7").Value = 0.0008036091882046079
$ws.Range("P7").Value = 0.000803609188204608
$ws.Range("Q7").Value = 3.022381720066666
$ws.Range("R7").Value = 27.2014354806
$ws.Range("S7").Value = 0.00007105948995880715
$ws.Range("T7").Value = 0.00007105948995880717
$ws.Range("I8").Value = 0.08842543241393927
$ws.Range("J8").Value = 0.08842543241393927
$ws.Range("M8").Value = 23.13337333333333
$ws.Range("N8").Value = 69.40011999999999
$ws.Range("O8").Value = 0.4750474795102416
$ws.Range("P8").Value = 0.4750474795102416
$ws.Range("Q8").Value = 1786.658041383586
$ws.Range("R8").Value = 16079.92237245228
$ws.Range("S8").Value = 0.04200627879284507
$ws.Range("T8").Value = 0.04200627879284507
$ws.Range("I9").Value = 0.08842543241393927
$ws.Range("J9").Value = 0.08842543241393927
$ws.Range("M9").Value = 13.33195766666667
$ws.Range("N9").Value = 39.995873
$ws.Range("O9").Value = 0.273773858884707
$ws.Range("P9").Value = 0.273773858884707
$ws.Range("Q9").Value = 1029.666059908926
$ws.Range("R9").Value = 9266.994539180338
$ws.Range("S9").Value = 0.02420857185551301
$ws.Range("T9").Value = 0.024208571855513
$ws.Range("G10").Value = 174.3107043333333
$ws.Range("H10").Value = 522.932113
$ws.Range("I10").Value = 0.199571617988009
$ws.Range("J10").Value = 0.199571617988009
$ws.Range("M10").Value = 12.19250666666667
$ws.Range("N10").Value = 36.57752
$ws.Range("O10").Value = 0.2503750524168468
$ws.Range("P10").Value = 0.2503750524168468
$ws.Range("Q10").Value = 2125.284424655529
$ws.Range("R10").Value = 19127.55982189976
$ws.Range("S10").Value = 0.04996775431466268
$ws.Range("T10").Value = 0.04996775431466267
$ws.Range("G11").Value = 174.3107043333333
$ws.Range("H11").Value = 522.932113
$ws.Range("I11").Value = 0.199571617988009
$ws.Range("J11").Value = 0.199571617988009
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.03913333333333333
$ws.Range("N11").Value = 0.1174
$ws.Range("O11").Value = 0.0008036091882046079
$ws.Range("P11").Value = 0.000803609188204608
$ws.Range("Q11").Value = 6.821358896244444
$ws.Range("R11").Value = 61.3922300662
$ws.Range("S11").Value = 0.0001603775859200241
$ws.Range("T11").Value = 0.0001603775859200241
$ws.Range("G12").Value = 174.3107043333333
$ws.Range("H12").Value = 522.932113
$ws.Range("I12").Value = 0.199571617988009
$ws.Range("J12").Value = 0.199571617988009
$ws.Range("M12").Value = 23.13337333333333
$ws.Range("N12").Value = 69.40011999999999
$ws.Range("O12").Value = 0.4750474795102416
$ws.Range("P12").Value = 0.4750474795102416
$ws.Range("Q12").Value = 4032.394599339283
$ws.Range("R12").Value = 36291.55139405355
$ws.Range("S12").Value = 0.09480599410698448
$ws.Range("T12").Value = 0.09480599410698447
$ws.Range("G13").Value = 174.3107043333333
$ws.Range("H13").Value = 522.932113
$ws.Range("I13").Value = 0.199571617988009
$ws.Range("J13").Value = 0.199571617988009
$ws.Range("M13").Value = 13.33195766666667
$ws.Range("N13").Value = 39.995873
$ws.Range("O13").Value = 0.273773858884707
$ws.Range("P13").Value = 0.273773858884707
$ws.Range("Q13").Value = 2323.90293101885
$ws.Range("R13").Value = 20915.12637916965
$ws.Range("S13").Value = 0.05463749198044184
$ws.Range("T13").Value = 0.05463749198044182
$ws.Range("G14").Value = 28.53474833333333
$ws.Range("H14").Value = 85.60424499999999
$ws.Range("I14").Value = 0.03266997236655063
$ws.Range("J14").Value = 0.03266997236655063
$ws.Range("M14").Value = 12.19250666666667
$ws.Range("N14").Value = 36.57752
$ws.Range("O14").Value = 0.2503750524168468
$ws.Range("P14").Value = 0.2503750524168468
$ws.Range("Q14").Value = 347.9101092858222
$ws.Range("R14").Value = 3131.1909835724
$ws.Range("S14").Value = 0.008179746043732049
$ws.Range("T14").Value = 0.008179746043732049
$ws.Range("G15").Value = 28.53474833333333
$ws.Range("H15").Value = 85.60424499999999
$ws.Range("I15").Value = 0.03266997236655063
$ws.Range("J15").Value = 0.03266997236655063
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 0.6666666666666666
$ws.Range("M15").Value = 0.03913333333333333
$ws.Range("N15").Value = 0.1174
$ws.Range("O15").Value = 0.0008036091882046079
$ws.Range("P15").Value = 0.000803609188204608
$ws.Range("Q15").Value = 1.116659818111111
$ws.Range("R15").Value = 10.049938363
$ws.Range("S15").Value = 0.00002625388997215072
$ws.Range("T15").Value = 0.00002625388997215073
$ws.Range("G16").Value = 28.53474833333333
$ws.Range("H16").Value = 85.60424499999999
$ws.Range("I16").Value = 0.03266997236655063
$ws.Range("J16").Value = 0.03266997236655063
$ws.Range("M16").Value = 23.13337333333333
$ws.Range("N16").Value = 69.40011999999999
$ws.Range("O16").Value = 0.4750474795102416
$ws.Range("P16").Value = 0.4750474795102416
$ws.Range("Q16").Value = 660.1049861677109
$ws.Range("R16").Value = 5940.944875509398
$ws.Range("S16").Value = 0.01551978802839912
$ws.Range("T16").Value = 0.01551978802839912
$ws.Range("G17").Value = 28.53474833333333
$ws.Range("H17").Value = 85.60424499999999
$ws.Range("I17").Value = 0.03266997236655063
$ws.Range("J17").Value = 0.03266997236655063
$ws.Range("M17").Value = 13.33195766666667
$ws.Range("N17").Value = 39.995873
$ws.Range("O17").Value = 0.273773858884707
$ws.Range("P17").Value = 0.273773858884707
$ws.Range("Q17").Value = 380.4240568089872
$ws.Range("R17").Value = 3423.816511280885
$ws.Range("S17").Value = 0.00894418440444731
$ws.Range("T17").Value = 0.008944184404447308
$ws.Range("G18").Value = 230.32901
$ws.Range("H18").Value = 690.98703
$ws.Range("I18").Value = 0.263708034289011
$ws.Range("J18").Value = 0.263708034289011
$ws.Range("M18").Value = 12.19250666666667
$ws.Range("N18").Value = 36.57752
$ws.Range("O18").Value = 0.2503750524168468
$ws.Range("P18").Value = 0.2503750524168468
$ws.Range("Q18").Value = 2808.287989951733
$ws.Range("R18").Value = 25274.5919095656
$ws.Range("S18").Value = 0.06602591290785474
$ws.Range("T18").Value = 0.06602591290785474
$ws.Range("G19").Value = 230.32901
$ws.Range("H19").Value = 690.98703
$ws.Range("I19").Value = 0.263708034289011
$ws.Range("J19").Value = 0.263708034289011
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 0.03913333333333333
$ws.Range("N19").Value = 0.1174
$ws.Range("O19").Value = 0.0008036091882046079
$ws.Range("P19").Value = 0.000803609188204608
$ws.Range("Q19").Value = 9.013541924666667
$ws.Range("R19").Value = 81.121877322
$ws.Range("S19").Value = 0.000211918199358025
$ws.Range("T19").Value = 0.000211918199358025
$ws.Range("G20").Value = 230.32901
$ws.Range("H20").Value = 690.98703
$ws.Range("I20").Value = 0.263708034289011
$ws.Range("J20").Value = 0.263708034289011
$ws.Range("M20").Value = 23.13337333333333
$ws.Range("N20").Value = 69.40011999999999
$ws.Range("O20").Value = 0.4750474795102416
$ws.Range("P20").Value = 0.4750474795102416
$ws.Range("Q20").Value = 5328.286977827066
$ws.Range("R20").Value = 47954.58280044359
$ws.Range("S20").Value = 0.125273837015595
$ws.Range("T20").Value = 0.125273837015595
$ws.Range("G21").Value = 230.32901
$ws.Range("H21").Value = 690.98703
$ws.Range("I21").Value = 0.263708034289011
$ws.Range("J21").Value = 0.263708034289011
$ws.Range("M21").Value = 13.33195766666667
$ws.Range("N21").Value = 39.995873
$ws.Range("O21").Value = 0.273773858884707
$ws.Range("P21").Value = 0.273773858884707
$ws.Range("Q21").Value = 3070.736610725244
$ws.Range("R21").Value = 27636.62949652719
$ws.Range("S21").Value = 0.07219636616620316
$ws.Range("T21").Value = 0.07219636616620315
$ws.Range("G22").Value = 151.5554656666667
$ws.Range("H22").Value = 454.666397
$ws.Range("I22").Value = 0.1735187154093718
$ws.Range("J22").Value = 0.1735187154093718
$ws.Range("M22").Value = 12.19250666666667
$ws.Range("N22").Value = 36.57752
$ws.Range("O22").Value = 0.2503750524168468
$ws.Range("P22").Value = 0.2503750524168468
$ws.Range("Q22").Value = 1847.841025510604
$ws.Range("R22").Value = 16630.56922959544
$ws.Range("S22").Value = 0.04344475746592538
$ws.Range("T22").Value = 0.04344475746592538
$ws.Range("G23").Value = 151.5554656666667
$ws.Range("H23").Value = 454.666397
$ws.Range("I23").Value = 0.1735187154093718
$ws.Range("J23").Value = 0.1735187154093718
$ws.Range("K23").Value = 2
$ws.Range("L23").Value = 0.6666666666666666
$ws.Range("M23").Value = 0.03913333333333333
$ws.Range("N23").Value = 0.1174
$ws.Range("O23").Value = 0.0008036091882046079
$ws.Range("P23").Value = 0.000803609188204608
$ws.Range("Q23").Value = 5.930870556422222
$ws.Range("R23").Value = 53.37783500779999
$ws.Range("S23").Value = 0.0001394412340284317
$ws.Range("T23").Value = 0.0001394412340284317
$ws.Range("G24").Value = 151.5554656666667
$ws.Range("H24").Value = 454.666397
$ws.Range("I24").Value = 0.1735187154093718
$ws.Range("J24").Value = 0.1735187154093718
$ws.Range("M24").Value = 23.13337333333333
$ws.Range("N24").Value = 69.40011999999999
$ws.Range("O24").Value = 0.4750474795102416
$ws.Range("P24").Value = 0.4750474795102416
$ws.Range("Q24").Value = 3505.989167974181
$ws.Range("R24").Value = 31553.90251176763
$ws.Range("S24").Value = 0.08242962840307699
$ws.Range("T24").Value = 0.08242962840307699
$ws.Range("G25").Value = 151.5554656666667
$ws.Range("H25").Value = 454.666397
$ws.Range("I25").Value = 0.1735187154093718
$ws.Range("J25").Value = 0.1735187154093718
$ws.Range("M25").Value = 13.33195766666667
$ws.Range("N25").Value = 39.995873
$ws.Range("O25").Value = 0.273773858884707
$ws.Range("P25").Value = 0.273773858884707
$ws.Range("Q25").Value = 2020.531052419954
$ws.Range("R25").Value = 18184.77947177958
$ws.Range("S25").Value = 0.04750488830634099
$ws.Range("T25").Value = 0.04750488830634098
